$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newD20 = @'
1) Users are to fill in all fields 
eg. Imageurl, name, email ,contact no, country origin, color, skin type, treat, ingredients, recommended usage and more…
2) For oil, base and milk ingredient input fill, users are to register input value by hitting enter key after keying in value 
3) Any fields that are not filled or not filled properly, validation error will be shown below each field upon clicking on add button
4) If no fields are filled and upon clicking on add button, validation errors will show for all fields
5) After filling all fields , click on add button located at the bottom of modal  to add to collection
'@

$newD16 = @'
1) Users are able to edit all fields in the collection 
eg. Imageurl, name, email ,contact no, country origin, color, skin type, treat, ingredients, recommended usage and more…
2) For oil, base and milk ingredient input fill, users are to register input value by hitting enter key after keying in value 
3) Any fields that are not filled or not filled properly, validation error will be shown below each field upon clicking on update button
4) If no fields are filled and upon clicking on update button, validation errors will show for all fields
5) After editing , click on update located at the bottom of modal  to update the changes
'@

$ws.Range("D20").Value = $newD20
$ws.Range("D16").Value = $newD16

$ws.Rows.Item(16).RowHeight = 249.6
$ws.Rows.Item(20).RowHeight = 246.6

$ws.Range("C16").Select()
